$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows with LDIF-style distinguished names
$ws.Range("A1").Value = "cn=shubham,ou=users,o=abc"
$ws.Range("B1").Value = "cn=group1,ou=group,o=abc"

$ws.Range("A2").Value = "cn=noob,ou=users,o=abc"
$ws.Range("B2").Value = "cn=group1,ou=group,o=abc"

# Add a new third row
$ws.Range("A3").Value = "cn=trump,ou=users,o=abc"
$ws.Range("B3").Value = "cn=group3,ou=server,o=abc"

# Resize columns to fit content, mirroring the bestFit column widths in the target file
$ws.Columns.Item(1).ColumnWidth = 23.053385416666668
$ws.Columns.Item(2).ColumnWidth = 22.166666666666668

# Set the active selection to B3, matching the target worksheet's selection
$ws.Range("B3").Select() | Out-Null
